$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Chemours El Dorado / El Dorado, AR
$ws.Cells.Item(2,1).Value = "Chemours El Dorado"
$ws.Cells.Item(2,2).Value = "El Dorado, AR"
$ws.Cells.Item(2,3).Value = 2
$ws.Cells.Item(2,4).Value = 2
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 12

# Row 3: San Dimas Plant / San Dimas, CA
$ws.Cells.Item(3,1).Value = "San Dimas Plant"
$ws.Cells.Item(3,2).Value = "San Dimas, CA"
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = 1
$ws.Cells.Item(3,5).Value = 4
$ws.Cells.Item(3,6).Value = 34

# Row 4: CF Industries Nitrogen LLC-Port Neal Nitrogen Complex / Sergeant Bluff, IA
$ws.Cells.Item(4,1).Value = "CF Industries Nitrogen LLC-Port Neal Nitrogen Complex"
$ws.Cells.Item(4,2).Value = "Sergeant Bluff, IA"
$ws.Cells.Item(4,3).Value = 2
$ws.Cells.Item(4,4).Value = 6
$ws.Cells.Item(4,5).Value = 7
$ws.Cells.Item(4,6).Value = 21

# Row 5: Linde - Whiting / East Chicago, IN
$ws.Cells.Item(5,1).Value = "Linde - Whiting"
$ws.Cells.Item(5,2).Value = "East Chicago, IN"
$ws.Cells.Item(5,3).Value = 5
$ws.Cells.Item(5,4).Value = 27
$ws.Cells.Item(5,5).Value = 35
$ws.Cells.Item(5,6).Value = 71

# Row 6: AIR PRODUCTS & CHEMICALS INC - Geismar SMR / Geismar, LA
$ws.Cells.Item(6,1).Value = "AIR PRODUCTS & CHEMICALS INC - Geismar SMR"
$ws.Cells.Item(6,2).Value = "Geismar, LA"
$ws.Cells.Item(6,3).Value = 3
$ws.Cells.Item(6,4).Value = 13
$ws.Cells.Item(6,5).Value = 18
$ws.Cells.Item(6,6).Value = 42

# Row 7: HC Manvel Inc / Manvel, TX
$ws.Cells.Item(7,1).Value = "HC Manvel Inc"
$ws.Cells.Item(7,2).Value = "Manvel, TX"
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(7,4).Value = 1
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 10

# Row 8: Air Products Port Arthur Facility / Port Arthur, TX
$ws.Cells.Item(8,1).Value = "Air Products Port Arthur Facility"
$ws.Cells.Item(8,2).Value = "Port Arthur, TX"
$ws.Cells.Item(8,3).Value = 2
$ws.Cells.Item(8,4).Value = 15
$ws.Cells.Item(8,5).Value = 15
$ws.Cells.Item(8,6).Value = 31

# Row 9: KSP CO2 Plant / Tad, WV
$ws.Cells.Item(9,1).Value = "KSP CO2 Plant"
$ws.Cells.Item(9,2).Value = "Tad, WV"
$ws.Cells.Item(9,3).Value = $null
$ws.Cells.Item(9,4).Value = $null
$ws.Cells.Item(9,5).Value = $null
$ws.Cells.Item(9,6).Value = 3

# Row 10: Linde Decatur / Decatur, AL
$ws.Cells.Item(10,1).Value = "Linde Decatur"
$ws.Cells.Item(10,2).Value = "Decatur, AL"
$ws.Cells.Item(10,3).Value = 3
$ws.Cells.Item(10,4).Value = 11
$ws.Cells.Item(10,5).Value = 23
$ws.Cells.Item(10,6).Value = 29

# Row 11: CALAMCO / Stockton, CA
$ws.Cells.Item(11,1).Value = "CALAMCO"
$ws.Cells.Item(11,2).Value = "Stockton, CA"
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 7
$ws.Cells.Item(11,5).Value = 14
$ws.Cells.Item(11,6).Value = 22

# Row 12: Diversified CPC International / Channahon, IL
$ws.Cells.Item(12,1).Value = "Diversified CPC International"
$ws.Cells.Item(12,2).Value = "Channahon, IL"
$ws.Cells.Item(12,3).Value = 5
$ws.Cells.Item(12,4).Value = 6
$ws.Cells.Item(12,5).Value = 9
$ws.Cells.Item(12,6).Value = 24

# Row 13: Aeropres-Sibley / Sibley, LA (new row, no counts)
$ws.Cells.Item(13,1).Value = "Aeropres-Sibley"
$ws.Cells.Item(13,2).Value = "Sibley, LA"
$ws.Cells.Item(13,3).Value = $null
$ws.Cells.Item(13,4).Value = $null
$ws.Cells.Item(13,5).Value = $null
$ws.Cells.Item(13,6).Value = $null

# Row 14: Chemours-Corpus Christi / Gregory, TX (new row)
$ws.Cells.Item(14,1).Value = "Chemours-Corpus Christi"
$ws.Cells.Item(14,2).Value = "Gregory, TX"
$ws.Cells.Item(14,3).Value = 2
$ws.Cells.Item(14,4).Value = 4
$ws.Cells.Item(14,5).Value = 6
$ws.Cells.Item(14,6).Value = 6

# Row 15: Honeywell-Geismar / Geismar, LA (new row)
$ws.Cells.Item(15,1).Value = "Honeywell-Geismar"
$ws.Cells.Item(15,2).Value = "Geismar, LA"
$ws.Cells.Item(15,3).Value = 4
$ws.Cells.Item(15,4).Value = 21
$ws.Cells.Item(15,5).Value = 31
$ws.Cells.Item(15,6).Value = 36
